# Applies the cryptos.xlsx price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.124.84"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "1.676.68"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.43"
$ws.Range("E5").Value = "  -2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5275"
$ws.Range("E6").Value = "  -4.17%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2675"
$ws.Range("E8").Value = "  -0.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06328"
$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.30"
$ws.Range("E10").Value = "  -3.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07617"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("D12").Value = "1.690.44"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.507"
$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5690"
$ws.Range("E14").Value = "  -2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008177"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.88"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").Value = "26.169.42"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.840"
$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.63"
$ws.Range("E20").Value = "  -2.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.17"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.205"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.94"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1254"
$ws.Range("E25").Value = "  -4.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.644"
$ws.Range("E26").Value = "  -3.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.81"
$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06436"
$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("E29").Value = "  -2.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.300"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.542"
$ws.Range("E31").Value = "  -1.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.536"
$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.661"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.013"
$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.420"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6021"
$ws.Range("E36").Value = "  -2.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.136"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01625"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").Value = "1.089.89"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8699"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.05"
$ws.Range("E43").Value = "  -0.70%  "

$ws.Range("D44").Value = "1.829.53"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.06"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05256"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.996"
$ws.Range("E49").Value = "  -2.48%  "

$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.952"
$ws.Range("E51").Value = "  -1.85%  "
